$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) relabeling ---
$ws.Range("A1").Value = "Population Group"
$ws.Range("B1").Value = "All Institutions"
$ws.Range("C1").Value = "Total (4 Year)"
$ws.Range("D1").Value = "Male (4 Year)"
$ws.Range("E1").Value = "Female (4 Year)"
$ws.Range("F1").Value = "Total (2 Year)"
$ws.Range("G1").Value = "Male (2 Year)"
$ws.Range("H1").Value = "Female (2 Year)"
$ws.Range("I1").Value = "Total (Less than 2 Year)"
$ws.Range("J1").Value = "Male (Less Than 2 Year)"
$ws.Range("K1").Value = "Female (Less Than 2 Year)"

# --- Column A category relabeling (rows 2-21) ---
$ws.Range("A2").Value = "Total Applicants"
$ws.Range("A3").Value = "Total Admitted"
$ws.Range("A4").Value = "Total Enrollment"
$ws.Range("A5").Value = "Total Full Time"
$ws.Range("A6").Value = "Total Part Time"
$ws.Range("A7").Value = "Applicants (Public University)"
$ws.Range("A8").Value = "Admitted (Public University)"
$ws.Range("A9").Value = "Enrollment (Public University)"
$ws.Range("A10").Value = "Full Time Students (Public University)"
$ws.Range("A11").Value = "Part Time Students (Public University)"
$ws.Range("A12").Value = "Applicants (Private University)"
$ws.Range("A13").Value = "Admitted (Private University)"
$ws.Range("A14").Value = "Enrollment (Private University)"
$ws.Range("A15").Value = "Full Time Students (Private University)"
$ws.Range("A16").Value = "Part Time Students (Private University)"
$ws.Range("A17").Value = "Applicant (Private University)"
$ws.Range("A18").Value = "Admitted  (Private For Profit University)"
$ws.Range("A19").Value = "Enrollment (Private For Profit University)"
$ws.Range("A20").Value = "Full Time (Private For Profit University)"
$ws.Range("A21").Value = "Part Time (Private For Profit University)"

# --- Update the active selection to match the saved view state ---
$ws.Range("J1").Select()
